$wb = $excel.ActiveWorkbook

$oldUuid = "f5529299-f690-432a-b166-858b9977eb6b"
$newUuid = "32273afd-843c-479e-9c38-f887f81264b8"
$newZhHash = "$newUuid.9879430444c3287319a5818de9e09eaa6875f231.zh-cn.xlf"
$newDeHash = "$newUuid.9879430444c3287319a5818de9e09eaa6875f231.de-de.xlf"

# --- Sheet "Overview" ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A2").Value = "$newUuid.md"
foreach ($hl in $ws1.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$newUuid.md"
    }
}
$ws1.Range("G2").Value = "2016-08-23 22:55:58"

# --- Sheet "zh-cn" ---
$ws2 = $wb.Worksheets.Item("zh-cn")
foreach ($hl in $ws2.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$newUuid.md"
    } elseif ($addr -eq '$I$2') {
        $hl.Delete()
    }
}
$ws2.Range("A2").Value = "$newUuid.md"
$ws2.Range("G2").Value = $newZhHash
$ws2.Range("H2").Value = "2016-08-23 22:55:53"
$ws2.Range("I2").Value = ""
$ws2.Range("I2").Style = "Normal"
$ws2.Range("J2").Value = ""
$ws2.Range("K2").Value = "0001-01-01 00:00:00"
$ws2.Columns.Item(9).ColumnWidth = 18.6506053379604
$ws2.Columns.Item(10).ColumnWidth = 21.7054770333426

# --- Sheet "de-de" ---
$ws3 = $wb.Worksheets.Item("de-de")
foreach ($hl in $ws3.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$newUuid.md"
    } elseif ($addr -eq '$I$2') {
        $hl.Delete()
    }
}
$ws3.Range("A2").Value = "$newUuid.md"
$ws3.Range("G2").Value = $newDeHash
$ws3.Range("H2").Value = "2016-08-23 22:55:58"
$ws3.Range("I2").Value = ""
$ws3.Range("I2").Style = "Normal"
$ws3.Range("J2").Value = ""
$ws3.Range("K2").Value = "0001-01-01 00:00:00"
$ws3.Columns.Item(9).ColumnWidth = 18.6506053379604
$ws3.Columns.Item(10).ColumnWidth = 21.7054770333426
